# CHU_YR_FIN.xlsx — "Doing Updates for Financials"
# Refreshed the financial figures pulled from the data source (values are
# uniformly restated ~2% higher across the board, plus a handful of cells
# that are now reported as "NA" instead of a stale numeric figure).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D8").Value = 40787400
$ws.Range("E8").Value = 40693600
$ws.Range("F8").Value = 41116800
$ws.Range("G8").Value = 42249500
$ws.Range("H8").Value = 43786600
$ws.Range("I8").Value = 36943100
$ws.Range("J8").Value = 31042500
$ws.Range("D9").Value = 18396800
$ws.Range("E9").Value = 19533100
$ws.Range("F9").Value = 18979700
$ws.Range("G9").Value = 19413800
$ws.Range("H9").Value = 17412600
$ws.Range("I9").Value = 14282500
$ws.Range("J9").Value = 11215000
$ws.Range("D10").Value = 22390600
$ws.Range("E10").Value = 21160400
$ws.Range("F10").Value = 22137100
$ws.Range("G10").Value = 22835700
$ws.Range("H10").Value = 26373900
$ws.Range("I10").Value = 22660600
$ws.Range("J10").Value = 19827400
$ws.Range("F14").Value = -1372200
$ws.Range("D15").Value = 11500600
$ws.Range("E15").Value = 11398600
$ws.Range("F15").Value = 11388700
$ws.Range("G15").Value = 10962700
$ws.Range("H15").Value = 10121000
$ws.Range("I15").Value = 9061500
$ws.Range("J15").Value = 8610900
$ws.Range("D17").Value = 40203700
$ws.Range("E17").Value = 40293900
$ws.Range("F17").Value = 38147200
$ws.Range("G17").Value = 39444100
$ws.Range("H17").Value = 41446600
$ws.Range("I17").Value = 35221300
$ws.Range("J17").Value = 30242400
$ws.Range("D18").Value = 583700
$ws.Range("E18").Value = 399700
$ws.Range("F18").Value = 2969700
$ws.Range("G18").Value = 2805400
$ws.Range("H18").Value = 2340000
$ws.Range("I18").Value = 1721900
$ws.Range("J18").Value = 800100
$ws.Range("D20").Value = 578400
$ws.Range("E20").Value = 445100
$ws.Range("F20").Value = -214000
$ws.Range("G20").Value = 203200
$ws.Range("H20").Value = 294100
$ws.Range("I20").Value = 155500
$ws.Range("J20").Value = 377100
$ws.Range("D21").Value = 12682800
$ws.Range("E21").Value = 12263300
$ws.Range("F21").Value = 14164300
$ws.Range("G21").Value = 13990500
$ws.Range("H21").Value = 12772800
$ws.Range("I21").Value = 10954700
$ws.Range("J21").Value = "NA"
$ws.Range("D22").Value = 777200
$ws.Range("E22").Value = 728400
$ws.Range("F22").Value = 672700
$ws.Range("G22").Value = 644200
$ws.Range("H22").Value = 598800
$ws.Range("I22").Value = 464400
$ws.Range("J22").Value = 346400
$ws.Range("D23").Value = 384800
$ws.Range("E23").Value = 116400
$ws.Range("F23").Value = 2082900
$ws.Range("G23").Value = 2364300
$ws.Range("H23").Value = 2035300
$ws.Range("I23").Value = 1413000
$ws.Range("J23").Value = 830800
$ws.Range("D24").Value = 110300
$ws.Range("E24").Value = 22900
$ws.Range("F24").Value = 515400
$ws.Range("G24").Value = 575200
$ws.Range("H24").Value = 490600
$ws.Range("I24").Value = 359900
$ws.Range("J24").Value = 203500
$ws.Range("D26").Value = 274600
$ws.Range("E26").Value = 93500
$ws.Range("F26").Value = 1567500
$ws.Range("G26").Value = 1789100
$ws.Range("H26").Value = 1544700
$ws.Range("I26").Value = 1053100
$ws.Range("J26").Value = 627300
$ws.Range("D27").Value = 271300
$ws.Range("E27").Value = 92800
$ws.Range("F27").Value = 1567500
$ws.Range("G27").Value = 1789100
$ws.Range("H27").Value = 1544700
$ws.Range("I27").Value = 1053100
$ws.Range("J27").Value = 627300
$ws.Range("D32").Value = -578400
$ws.Range("E32").Value = -445100
$ws.Range("F32").Value = 214000
$ws.Range("G32").Value = -203200
$ws.Range("H32").Value = -294100
$ws.Range("I32").Value = -155500
$ws.Range("J32").Value = -377100
$ws.Range("D33").Value = 271300
$ws.Range("E33").Value = 92800
$ws.Range("F33").Value = 1567500
$ws.Range("G33").Value = 1789100
$ws.Range("H33").Value = 1544700
$ws.Range("I33").Value = 1053100
$ws.Range("J33").Value = 627300
$ws.Range("D35").Value = 271300
$ws.Range("E35").Value = 92800
$ws.Range("F35").Value = 1567500
$ws.Range("G35").Value = 1789100
$ws.Range("H35").Value = 1544700
$ws.Range("I35").Value = 1053100
$ws.Range("J35").Value = 627300
$ws.Range("D41").Value = 4873200
$ws.Range("E41").Value = 3507400
$ws.Range("F41").Value = 3228700
$ws.Range("G41").Value = 3756000
$ws.Range("H41").Value = 3191700
$ws.Range("I41").Value = 2708500
$ws.Range("J41").Value = 2241900
$ws.Range("D42").Value = 843900
$ws.Range("E42").Value = 278600
$ws.Range("F42").Value = 45700
$ws.Range("G42").Value = 8300
$ws.Range("H42").Value = 8000
$ws.Range("J42").Value = 42600
$ws.Range("D43").Value = 4352700
$ws.Range("E43").Value = 7193700
$ws.Range("F43").Value = 3747400
$ws.Range("G43").Value = 2916100
$ws.Range("H43").Value = 2806700
$ws.Range("I43").Value = 2153300
$ws.Range("J43").Value = 1872200
$ws.Range("D44").Value = 332300
$ws.Range("E44").Value = 360800
$ws.Range("F44").Value = 585600
$ws.Range("G44").Value = 649700
$ws.Range("H44").Value = 821600
$ws.Range("I44").Value = 861200
$ws.Range("J44").Value = 690300
$ws.Range("D45").Value = 984300
$ws.Range("E45").Value = 861500
$ws.Range("F45").Value = 803000
$ws.Range("G45").Value = 1066000
$ws.Range("H45").Value = 920400
$ws.Range("I45").Value = 1425900
$ws.Range("J45").Value = 911800
$ws.Range("D46").Value = 11386300
$ws.Range("E46").Value = 12202000
$ws.Range("F46").Value = 8410400
$ws.Range("G46").Value = 8396100
$ws.Range("H46").Value = 7748500
$ws.Range("I46").Value = 7149500
$ws.Range("J46").Value = 5758800
$ws.Range("D47").Value = 5919600
$ws.Range("E47").Value = 5602300
$ws.Range("F47").Value = 8333100
$ws.Range("G47").Value = 1326600
$ws.Range("H47").Value = 964200
$ws.Range("I47").Value = 826200
$ws.Range("J47").Value = 1031600
$ws.Range("D48").Value = 61827000
$ws.Range("E48").Value = 66950000
$ws.Range("F48").Value = 67471800
$ws.Range("G48").Value = 65051200
$ws.Range("H48").Value = 64057500
$ws.Range("I48").Value = 65092300
$ws.Range("J48").Value = 57808100
$ws.Range("D49").Value = 2042000
$ws.Range("E49").Value = 2061600
$ws.Range("F49").Value = 2001300
$ws.Range("G49").Value = 1761800
$ws.Range("H49").Value = 1616500
$ws.Range("I49").Value = 1355700
$ws.Range("J49").Value = 1228400
$ws.Range("D52").Value = 3713100
$ws.Range("E52").Value = 4330800
$ws.Range("F52").Value = 4364900
$ws.Range("G52").Value = 4358400
$ws.Range("H52").Value = 4147600
$ws.Range("I52").Value = 2174200
$ws.Range("J52").Value = 1882700
$ws.Range("D54").Value = 84888000
$ws.Range("E54").Value = 91146600
$ws.Range("F54").Value = 90581400
$ws.Range("G54").Value = 80894100
$ws.Range("H54").Value = 78534300
$ws.Range("I54").Value = 76598000
$ws.Range("J54").Value = 67709500
$ws.Range("D57").Value = 12917800
$ws.Range("E57").Value = 16436000
$ws.Range("F57").Value = 20220400
$ws.Range("G57").Value = 14995500
$ws.Range("H57").Value = 11610600
$ws.Range("I57").Value = 12954400
$ws.Range("J57").Value = 11467800
$ws.Range("D58").Value = 7468300
$ws.Range("E58").Value = 19987100
$ws.Range("F58").Value = 15828500
$ws.Range("G58").Value = 16749800
$ws.Range("H58").Value = 19214600
$ws.Range("I58").Value = 20664600
$ws.Range("J58").Value = 10443900
$ws.Range("D59").Value = 15621500
$ws.Range("E59").Value = 14430300
$ws.Range("F59").Value = 13827800
$ws.Range("G59").Value = 11578500
$ws.Range("H59").Value = 12991200
$ws.Range("I59").Value = 11248300
$ws.Range("J59").Value = 9837200
$ws.Range("D60").Value = 36007500
$ws.Range("E60").Value = 50853400
$ws.Range("F60").Value = 49876700
$ws.Range("G60").Value = 43323800
$ws.Range("H60").Value = 43816400
$ws.Range("I60").Value = 44867300
$ws.Range("J60").Value = 31748900
$ws.Range("D61").Value = 3218300
$ws.Range("E61").Value = 6022300
$ws.Range("F61").Value = 6076500
$ws.Range("G61").Value = 3561500
$ws.Range("H61").Value = 2038900
$ws.Range("I61").Value = 376400
$ws.Range("J61").Value = 5120400
$ws.Range("D62").Value = 494100
$ws.Range("E62").Value = 480600
$ws.Range("F62").Value = 313400
$ws.Range("G62").Value = 239400
$ws.Range("H62").Value = 192200
$ws.Range("I62").Value = 261600
$ws.Range("J62").Value = 282900
$ws.Range("D66").Value = 39763900
$ws.Range("E66").Value = 57397100
$ws.Range("F66").Value = 56266700
$ws.Range("G66").Value = 47124800
$ws.Range("H66").Value = 46047500
$ws.Range("I66").Value = 45505300
$ws.Range("J66").Value = 37152200
$ws.Range("D72").Value = 8457600
$ws.Range("E72").Value = 8198300
$ws.Range("F72").Value = 8685000
$ws.Range("G72").Value = 7727100
$ws.Range("H72").Value = 6139300
$ws.Range("I72").Value = 5483900
$ws.Range("J72").Value = 4780700
$ws.Range("D76").Value = 45124100
$ws.Range("E76").Value = 33749500
$ws.Range("F76").Value = 34314800
$ws.Range("G76").Value = 33769400
$ws.Range("H76").Value = 32486800
$ws.Range("I76").Value = 31092600
$ws.Range("J76").Value = 30557300
$ws.Range("D81").Value = 271300
$ws.Range("E81").Value = 92800
$ws.Range("F81").Value = 1567500
$ws.Range("G81").Value = 1789100
$ws.Range("H81").Value = 1544700
$ws.Range("I81").Value = 1053100
$ws.Range("J81").Value = 627300
$ws.Range("D83").Value = 11500600
$ws.Range("E83").Value = 11398600
$ws.Range("F83").Value = 11388700
$ws.Range("G83").Value = 10962700
$ws.Range("H83").Value = 10121000
$ws.Range("I83").Value = 9061500
$ws.Range("J83").Value = "NA"
$ws.Range("D89").Value = 12622900
$ws.Range("E89").Value = 11070300
$ws.Range("F89").Value = 12511100
$ws.Range("G89").Value = 13074000
$ws.Range("H89").Value = 11647500
$ws.Range("I89").Value = 10480700
$ws.Range("J89").Value = 9867900
$ws.Range("D91").Value = -9749500
$ws.Range("E91").Value = -15195000
$ws.Range("F91").Value = -13803200
$ws.Range("G91").Value = -10892300
$ws.Range("H91").Value = -11695900
$ws.Range("I91").Value = -12879500
$ws.Range("J91").Value = -11555400
$ws.Range("D94").Value = -7025100
$ws.Range("E94").Value = -14210100
$ws.Range("F94").Value = -13557800
$ws.Range("G94").Value = -11178100
$ws.Range("H94").Value = -11443900
$ws.Range("I94").Value = -14763800
$ws.Range("J94").Value = "NA"
$ws.Range("E96").Value = -604200
$ws.Range("F96").Value = -689100
$ws.Range("G96").Value = -545700
$ws.Range("H96").Value = -398600
$ws.Range("I96").Value = -338800
$ws.Range("J96").Value = -307200
$ws.Range("D100").Value = -4216900
$ws.Range("E100").Value = 3395200
$ws.Range("F100").Value = 508600
$ws.Range("G100").Value = -1331700
$ws.Range("H100").Value = 285800
$ws.Range("I100").Value = 4749700
$ws.Range("J100").Value = "NA"
$ws.Range("D101").Value = -15000
$ws.Range("E101").Value = 23300
$ws.Range("F101").Value = 10800
$ws.Range("H101").Value = -6200
$ws.Range("D102").Value = 1365800
$ws.Range("E102").Value = 278700
$ws.Range("F102").Value = -527300
$ws.Range("G102").Value = 564300
$ws.Range("H102").Value = 483200
$ws.Range("I102").Value = 466600
$ws.Range("J102").Value = -1111700
